# Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) across the resume bullets.
#
# Strategy: for each target paragraph, walk left-to-right through the list
# of metric substrings, each time restricting the Find to the remaining
# tail of the paragraph (from the end of the previous match onward) so
# that repeated/identical-looking tokens elsewhere in the document are
# never touched. Find.Execute collapses the supplied range to the match,
# and Word automatically splits the run there, so setting Font.Bold /
# Font.Color on that collapsed range only affects the matched substring -
# the surrounding plain text runs are left alone.

function Highlight-MetricsInParagraph($ParaIndex, $Phrases) {
    $para = $d.Paragraphs.Item($ParaIndex)
    $paraEnd = $para.Range.End
    $cursor = $para.Range.Start

    foreach ($phrase in $Phrases) {
        $searchRange = $d.Range($cursor, $paraEnd)
        $found = $searchRange.Find.Execute($phrase, $true, $false, $false, $false, $false, `
                                            $true, 1, $false, "", 0)
        if ($found) {
            $searchRange.Font.Bold = 1
            $searchRange.Font.Color = 5258796   # wdColor for RGB 2C3E50
            $cursor = $searchRange.End
        }
    }
}

$d = $word.ActiveDocument
$pm = [char]0x00B1   # '±' (U+00B1), built this way to avoid source-encoding issues

# Paragraph 10: "Discovered systematic race coding errors ... from 23% to 64%"
Highlight-MetricsInParagraph 10 @("23%", "64%")

# Paragraph 12: "Utilized advanced sampling methods ... ±4.2% to ±2.1% ...
#                increasing voter turnout prediction accuracy from 71% to 87% ..."
Highlight-MetricsInParagraph 12 @(($pm + "4.2%"), ($pm + "2.1%"), "71%", "87%")

# Paragraph 13: "Trigonometric algorithm ... reduced mapping costs by 73.5%,
#                saving campaigns and organizations $4.7M ..."
Highlight-MetricsInParagraph 13 @("73.5%", "`$4.7M")

# Paragraph 14: "Built real-time FEC analysis ... valued over $2 trillion"
Highlight-MetricsInParagraph 14 @("`$2")

# Paragraph 34: "Modernized legacy ETL processes ... reducing processing time by 57%"
Highlight-MetricsInParagraph 34 @("57%")

# Paragraph 50: "Predictive excellence: ... margin of error from ±4.2% to ±2.1%"
Highlight-MetricsInParagraph 50 @(($pm + "4.2%"), ($pm + "2.1%"))

# Paragraph 51: "Increased voter turnout prediction accuracy from 71% to 87%"
Highlight-MetricsInParagraph 51 @("71%", "87%")

# Paragraph 52: "Methodological advancement: Improved segmentation accuracy 34%
#                and survey incidence 28%"
Highlight-MetricsInParagraph 52 @("34%", "28%")
